$wb = $excel.ActiveWorkbook

function Set-HitsData {
    param($SheetName, $BValues, $DValues)
    $ws = $wb.Worksheets.Item($SheetName)
    for ($i = 0; $i -lt $BValues.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $BValues[$i]
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $DValues[$i]
    }
}

# "Total Hits" sheet
Set-HitsData "Total Hits" @(1959, 3917, 5858, 7773, 9741) @("52.77%", "52.76%", "52.60%", "52.35%", "52.48%")

# "Hits_entity" sheet
Set-HitsData "Hits_entity" @(1142, 2269, 3401, 4515, 5658) @("51.65%", "51.31%", "51.27%", "51.05%", "51.18%")

# "Hits_boolean" sheet
Set-HitsData "Hits_boolean" @(300, 605, 897, 1198, 1499) @("52.36%", "52.79%", "52.18%", "52.27%", "52.32%")
